$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("D1").Value = 1298
$ws.Range("F1").Value = 1298
$ws.Range("H1").Value = 1298
